$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 40, shifting existing rows 40.. down by one.
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new observation.
$ws.Range("A40").Value = 9
$ws.Range("B40").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C40").Value = "Metropolitana"
$ws.Range("D40").Value = 44930
$ws.Range("E40").Value = 13
$ws.Range("F40").Value = 100112029
$ws.Range("G40").Value = "Orégano"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 16
$ws.Range("K40").Value = 20000
$ws.Range("L40").Value = 20000
$ws.Range("M40").Value = 20000
$ws.Range("N40").Value = "$/docena de atados"
$ws.Range("O40").Value = "Región Metropolitana"
$ws.Range("P40").Value = 6667
$ws.Range("Q40").Value = 3
$ws.Range("R40").Value = "Hortaliza"
